$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 113.875
$ws.Range("I5").Value = 80.14286
$ws.Range("K5").Value = 80.14286
$ws.Range("M5").Value = 34.85714
$ws.Range("H18").Value = 5067
$ws.Range("I18").Value = 5067
$ws.Range("K18").Value = 5067
$ws.Range("M18").Value = -4783
$ws.Range("H98").Value = 874.1556
$ws.Range("I98").Value = 775.27905
$ws.Range("K98").Value = 775.27905
$ws.Range("M98").Value = 722.72095
$ws.Range("H105").Value = 49975
$ws.Range("J105").Value = 49975
$ws.Range("L105").Value = 49975
$ws.Range("N105").Value = -56963
$ws.Range("H122").Value = 874.1556
$ws.Range("I122").Value = 775.27905
$ws.Range("K122").Value = 2325.83715
$ws.Range("M122").Value = 124.1628500000002
$ws.Range("H127").Value = 42174.73
$ws.Range("J127").Value = 6146.75
$ws.Range("L127").Value = 18440.25
$ws.Range("N127").Value = -28360.25
$ws.Range("H137").Value = 13893485
$ws.Range("I137").Value = 23811972
$ws.Range("J137").Value = 7603.8667
$ws.Range("K137").Value = 71435916
$ws.Range("L137").Value = 22811.6001
$ws.Range("M137").Value = -71433366
$ws.Range("N137").Value = -27911.6001
$ws.Range("H138").Value = 4809.067
$ws.Range("I138").Value = 2991.5454
$ws.Range("J138").Value = 5397.0884
$ws.Range("K138").Value = 8974.636200000001
$ws.Range("L138").Value = 16191.2652
$ws.Range("M138").Value = -3834.636200000001
$ws.Range("N138").Value = -26471.2652
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2580.875
$ws.Range("I2").Value = 2580.875
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2580.875
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2467.875
$ws.Range("N2").ClearContents()
$ws.Range("H32").Value = 169750.97
$ws.Range("I32").Value = 201761.02
$ws.Range("K32").Value = 201761.02
$ws.Range("M32").Value = -201474.02
$ws.Range("H61").Value = 18186342
$ws.Range("I61").Value = 5346.25
$ws.Range("K61").Value = 5346.25
$ws.Range("M61").Value = -5134.25
$ws.Range("H88").Value = 2517.7273
$ws.Range("I88").Value = 2149
$ws.Range("J88").Value = 2728.4285
$ws.Range("K88").Value = 2149
$ws.Range("L88").Value = 2728.4285
$ws.Range("M88").Value = -1743
$ws.Range("N88").Value = -3540.4285
$ws.Range("H91").Value = 2517.7273
$ws.Range("I91").Value = 2149
$ws.Range("J91").Value = 2728.4285
$ws.Range("K91").Value = 2149
$ws.Range("L91").Value = 2728.4285
$ws.Range("M91").Value = -745
$ws.Range("N91").Value = -5536.4285
$ws.Range("H116").Value = 2580.875
$ws.Range("I116").Value = 2580.875
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2580.875
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -286.875
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 895068.2
$ws.Range("I132").Value = 1668543.1
$ws.Range("J132").Value = 2597.077
$ws.Range("K132").Value = 5005629.300000001
$ws.Range("L132").Value = 7791.231000000001
$ws.Range("M132").Value = -5003099.300000001
$ws.Range("N132").Value = -12851.231
$ws.Range("H136").Value = 18186342
$ws.Range("I136").Value = 5346.25
$ws.Range("K136").Value = 16038.75
$ws.Range("M136").Value = -13488.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2580.875
$ws.Range("I3").Value = 2580.875
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2580.875
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2466.875
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 4199.5
$ws.Range("I5").Value = 7398
$ws.Range("J5").Value = 1001
$ws.Range("K5").Value = 7398
$ws.Range("L5").Value = 1001
$ws.Range("M5").Value = -7285
$ws.Range("N5").Value = -1227
$ws.Range("H7").Value = 627537.25
$ws.Range("I7").Value = 835083
$ws.Range("K7").Value = 835083
$ws.Range("M7").Value = -834970
$ws.Range("H107").Value = 2666.6667
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = -80
$ws.Range("N107").Value = -7840
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5557637
$ws.Range("I31").Value = 6174824.5
$ws.Range("K31").Value = 6174824.5
$ws.Range("M31").Value = -6174529.5
$ws.Range("H34").Value = 5557637
$ws.Range("I34").Value = 6174824.5
$ws.Range("K34").Value = 6174824.5
$ws.Range("M34").Value = -6174622.5
$ws.Range("H52").Value = 49749.5
$ws.Range("I52").Value = 49749.5
$ws.Range("K52").Value = 49749.5
$ws.Range("M52").Value = -49455.5
$ws.Range("H107").Value = 568
$ws.Range("I107").Value = 439
$ws.Range("K107").Value = 439
$ws.Range("M107").Value = 1481
$ws.Range("H132").Value = 4685.05
$ws.Range("I132").Value = 5825.5
$ws.Range("J132").Value = 2974.375
$ws.Range("K132").Value = 17476.5
$ws.Range("L132").Value = 8923.125
$ws.Range("M132").Value = -14946.5
$ws.Range("N132").Value = -13983.125
$ws.Range("H134").Value = 2515.889
$ws.Range("I134").Value = 2028.24
$ws.Range("K134").Value = 6084.72
$ws.Range("M134").Value = -3549.72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6524.778
$ws.Range("J131").Value = 10145.4
$ws.Range("L131").Value = 30436.2
$ws.Range("N131").Value = -40516.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 20481
$ws.Range("I18").Value = 20481
$ws.Range("K18").Value = 20481
$ws.Range("M18").Value = -20188
$ws.Range("H21").Value = 181076.42
$ws.Range("I21").Value = 5494.6665
$ws.Range("K21").Value = 5494.6665
$ws.Range("M21").Value = -5321.6665
$ws.Range("H30").Value = 181076.42
$ws.Range("I30").Value = 5494.6665
$ws.Range("K30").Value = 5494.6665
$ws.Range("M30").Value = -5389.6665
$ws.Range("H97").Value = 980.1053000000001
$ws.Range("I97").Value = 891
$ws.Range("J97").Value = 1737.5
$ws.Range("K97").Value = 891
$ws.Range("L97").Value = 1737.5
$ws.Range("M97").Value = -395
$ws.Range("N97").Value = -2729.5
$ws.Range("H107").Value = 5979.222
$ws.Range("I107").Value = 7135.533
$ws.Range("J107").Value = 197.66667
$ws.Range("K107").Value = 7135.533
$ws.Range("L107").Value = 197.66667
$ws.Range("M107").Value = -5215.533
$ws.Range("N107").Value = -4037.66667
$ws.Range("H126").Value = 27602.143
$ws.Range("I126").Value = 43828.75
$ws.Range("K126").Value = 131486.25
$ws.Range("M126").Value = -129016.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 4987.6665
$ws.Range("I13").Value = 4987
$ws.Range("K13").Value = 4987
$ws.Range("M13").Value = -4847
$ws.Range("H23").Value = 9811
$ws.Range("I23").Value = 9802.556
$ws.Range("K23").Value = 9802.556
$ws.Range("M23").Value = -9572.556
$ws.Range("H122").Value = 8211
$ws.Range("I122").Value = 5948
$ws.Range("K122").Value = 17844
$ws.Range("M122").Value = -15394
$ws.Range("H132").Value = 4778788
$ws.Range("I132").Value = 9552720
$ws.Range("K132").Value = 28658160
$ws.Range("M132").Value = -28655630
$ws.Range("H136").Value = 7151884
$ws.Range("I136").Value = 6585226
$ws.Range("K136").Value = 19755678
$ws.Range("M136").Value = -19753128
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1999996.6
$ws.Range("I4").Value = 1999996.6
$ws.Range("K4").Value = 1999996.6
$ws.Range("M4").Value = -1999883.6
$ws.Range("H81").Value = 5752.067
$ws.Range("I81").Value = 2570.5715
$ws.Range("K81").Value = 5141.143
$ws.Range("M81").Value = -4080.143
$ws.Range("H84").Value = 5752.067
$ws.Range("I84").Value = 2570.5715
$ws.Range("K84").Value = 25705.715
$ws.Range("M84").Value = -20401.715
$ws.Range("H132").Value = 4632044
$ws.Range("I132").Value = 4764262
$ws.Range("K132").Value = 14292786
$ws.Range("M132").Value = -14290256
$ws.Range("H136").Value = 14807896
$ws.Range("I136").Value = 3108524.2
$ws.Range("K136").Value = 9325572.600000001
$ws.Range("M136").Value = -9323022.600000001
